$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove hyperlinks from B1 and B2 (and their underline/blue styling)
$ws.Range("B1").Hyperlinks.Delete()
$ws.Range("B2").Hyperlinks.Delete()

# Update the login detail values
$ws.Range("A1").Value = "testingblz123"
$ws.Range("B1").Value = "Theend@1"
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = ""

# Adjust column widths
$ws.Columns.Item(1).ColumnWidth = 21.85546875
$ws.Columns.Item(2).ColumnWidth = 19.42578125

# Apply a shared font style (JetBrains Mono, green, vertical center) to the data rows
$dataRange = $ws.Range("A1:B2")
$dataRange.Font.Name = "JetBrains Mono"
$dataRange.Font.Size = 9.8
$dataRange.Font.Color = 5867370
$dataRange.VerticalAlignment = -4108

# Move active selection to B1
$ws.Range("B1").Select()

# Configure page setup for printing
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
